$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.279703497886658
$ws.Range("B1").Value = 2.165440559387207
$ws.Range("C1").Value = 4.724639415740967
$ws.Range("D1").Value = 3.137783050537109
$ws.Range("E1").Value = 1.363385558128357
